# Practice Excel, using Solver
# Re-scope the LHS/RHS Solver constraint ranges into per-row single-cell
# references (splitting the old 2-row ranges into individual constraints),
# bump a couple of Solver engine flags, and refresh the worksheet data /
# objective comparison block to reflect a new Solver run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Solver defined names -------------------------------------------------
# Split solver_lhs1 (Sheet1!$B$15:$B$16) into two single-cell constraints,
# and move the old solver_lhs2 (Sheet1!$E$5:$E$10) down to solver_lhs3.
$wb.Names.Item("solver_lhs1").RefersTo = "=Sheet1!`$B`$15"
$wb.Names.Item("solver_lhs2").RefersTo = "=Sheet1!`$B`$16"
$ws.Names.Add("solver_lhs3", "=Sheet1!`$E`$5:`$E`$10")
$wb.Names.Item("solver_lhs3").Visible = $false

# Same split for the matching RHS names.
$wb.Names.Item("solver_rhs1").RefersTo = "=Sheet1!`$D`$15"
$wb.Names.Item("solver_rhs2").RefersTo = "=Sheet1!`$D`$16"
$ws.Names.Add("solver_rhs3", "=Sheet1!`$D`$5:`$D`$10")
$wb.Names.Item("solver_rhs3").Visible = $false

# New relational-operator flag for the third constraint group.
$ws.Names.Add("solver_rel3", "=1")
$wb.Names.Item("solver_rel3").Visible = $false

# Constraint count went from 2 groups to 3.
$wb.Names.Item("solver_num").RefersTo = "=3"

# New Solver engine options that appear alongside the extra constraint.
$ws.Names.Add("solver_est", "=1")
$wb.Names.Item("solver_est").Visible = $false
$ws.Names.Add("solver_nwt", "=1")
$wb.Names.Item("solver_nwt").Visible = $false

# Solver model version bump.
$wb.Names.Item("solver_ver").RefersTo = "=3"

# --- Worksheet data: new Solver solution ----------------------------------
$ws.Range("D5").Value = 90
$ws.Range("E5").Value = 90
$ws.Range("E8").Value = 1
$ws.Range("E10").Value = 16

# Objective comparison block (row 12): keep the live SUMPRODUCT in B12,
# record the prior objective in C12, the new objective in D12, and show
# the delta between them in E12.
$ws.Range("C12").Value = 120514
$ws.Range("D12").Value = 120654
$ws.Range("E12").Formula = "=C12-D12"

# --- Selection -------------------------------------------------------------
$ws.Range("E12").Select() | Out-Null
